# "Finally, there is an output. Output can be sorted by belt rank or priority rating"
#
# The schedule columns ("1st Class"/"2nd Class"/"3rd Class") are collapsed down
# to just two abbreviated class-slot columns, the now-unused 3rd Class column
# is cleared out, and a batch of new student rows (8-14) is appended below the
# existing roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "3rd Class" header is no longer needed - clear it but keep its style.
$ws.Range("I1").Value = $null

# Re-label the remaining two class-slot columns for the existing students and
# drop whatever used to live in the (now retired) "3rd Class" column.
$ws.Range("G2").Value = "Mon-Beg"
$ws.Range("H2").Value = "Wed-Beg"
$ws.Range("I2").Value = $null

$ws.Range("G3").Value = "Tue-Beg"
$ws.Range("H3").Value = "Thu-Beg"
$ws.Range("I3").Value = $null

$ws.Range("G4").Value = "Mon-Adv"
$ws.Range("H4").Value = "Wed-Adv"
$ws.Range("I4").Value = $null

$ws.Range("G5").Value = "Mon-Adv"
$ws.Range("H5").Value = "Fri-Adv"
$ws.Range("I5").Value = $null

$ws.Range("G6").Value = "Tue-Int"
$ws.Range("H6").Value = "Fri-Int"
$ws.Range("I6").Value = $null

$ws.Range("G7").Value = "Mon-Beg"
$ws.Range("H7").Value = "Wed-Beg"
$ws.Range("I7").Value = $null

# Append the newly-registered students.
$newStudents = @(
    @("Joe",    "Mama",     "Red",          "B", "N/A", "1 May, 2022", "Mon-Adv", "Tue-Adv"),
    @("Jody",   "Cox",      "Green",        "B", "N/A", "1 May, 2022", "Mon-Int", "Sat-InA"),
    @("Ben",    "Harmin",   "Black Stripe", "B", "N/A", "1 May, 2022", "Mon-Adv", "Sat-InA"),
    @("Julie",  "Summers",  "Black Stripe", "D", "N/A", "1 Jun, 2022", "Mon-Adv", "Sat-InA"),
    @("Lady",   "Fingers",  "Orange",       "B", "N/A", "1 May, 2022", "Tue-Beg", "Thu-Beg"),
    @("Fat",    "Buddha",   "Blue",         "B", "N/A", "1 May, 2022", "Tue-Int", "Fri-Int"),
    @("Evangy", "Bush",     "Purple",       "A", "N/A", "1 Apr, 2022", "Mon-Adv", "Fri-Adv")
)

$rowNum = 8
foreach ($student in $newStudents) {
    $ws.Range("A$rowNum").Value = $student[0]
    $ws.Range("B$rowNum").Value = $student[1]
    $ws.Range("C$rowNum").Value = $student[2]
    $ws.Range("D$rowNum").Value = $student[3]
    $ws.Range("E$rowNum").Value = $student[4]
    $ws.Range("F$rowNum").Value = $student[5]
    $ws.Range("G$rowNum").Value = $student[6]
    $ws.Range("H$rowNum").Value = $student[7]
    $rowNum++
}

# The new "Belt" column needs its own (narrower) custom width.
$ws.Columns.Item(3).ColumnWidth = 10.83

# Leave the selection where the author's cursor ended up after entering data.
$ws.Range("A15").Select()
